$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Depot #0 totals (row 2)
$ws.Range("D2").Value = "Depot Total Demand: 15.0"
$ws.Range("E2").Value = "Depot Total Distance: 4.352055072784424"

# Update Truck #0 max distance (row 3)
$ws.Range("B3").Value = "Truck Max Distance: 5.0"

# Row 6 becomes Truck #1's header row (was previously an "empty" truck placeholder)
$ws.Range("A6").Value = "Truck #1"
$ws.Range("B6").Value = "Truck Max Distance: 5.0"
$ws.Range("C6").Value = "Truck Node Count: 1"
$ws.Range("D6").Value = "Truck Total Demand: 5.0"
$ws.Range("E6").Value = "Truck Total Distance: 2.5632007122039795"

# Insert a new shipment row (row 8) with numeric data, shifting old row 8 (Depot #1) down to row 9
$ws.Rows("8").Insert()

$ws.Range("A8").Value = 17.0
$ws.Range("B8").Value = 5.0
$ws.Range("C8").Value = 11.399999618530273
$ws.Range("D8").Value = 12.0
